# Apply the "Add files via upload" edit to the Test Scenarios sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in the new test scenario row (row 7) - previously blank separator/total row.
# Set D7 before B7 so the shared-string table gets the two new strings in the
# same order as the target workbook (index 13 = D7 text, index 14 = B7 text).
$ws.Range("D7").Value = "Validte the functionality of the Dropdown Header > Your account."
$ws.Range("C7").Value = "FRS"
$ws.Range("B7").Value = "(TS_004)`nHeader_Your account"

# Match the taller row height used by the other data rows.
$ws.Range("B7").EntireRow.RowHeight = 30

# Remove the old SUM formula/total that used to live in F7.
$ws.Range("F7").ClearContents()

# Update the "Number of Test Cases" total for the Forget Password row.
$ws.Range("F6").Value = 27

# Move the active selection to E8, as recorded in the saved view state.
$ws.Range("E8").Select()
